# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de) and for each of the two data rows,
# the localized file has now been handed back and found in sync with the
# en-US source, so:
#   - the overall Status (column C) becomes "Handed back: in sync with en-US"
#   - the "Latest Target File" (F) and "Latest Handback File" (G) columns
#     get populated (pointing at the same md / xlf files already linked
#     from columns A / D), each as a hyperlink exactly like the existing
#     A/B/D hyperlink cells
#   - the "Latest Handback DateTime" (H) gets stamped with the handback time
#
# Helper: find the Hyperlink object already attached to a given cell
# address (e.g. "A2") on a worksheet, so the new F/G hyperlinks can reuse
# the very same target URL / display text.
function Get-HyperlinkAt($ws, $addr) {
    $target = $ws.Range($addr).Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            return $hl
        }
    }
    return $null
}

function Add-MirrorHyperlink($ws, $fromAddr, $toAddr) {
    $src = Get-HyperlinkAt $ws $fromAddr
    if ($src -ne $null) {
        $null = $ws.Hyperlinks.Add($ws.Range($toAddr), $src.Address, [Type]::Missing, [Type]::Missing, $src.TextToDisplay)
    }
}

$wb = $excel.ActiveWorkbook

# Status text changes everywhere it appears (Overview + both locale sheets).
foreach ($ws in $wb.Worksheets) {
    $null = $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# zh-cn: row 2 (25c8b8b1...) and row 3 (c582a784...)
Add-MirrorHyperlink $zh "A2" "F2"
Add-MirrorHyperlink $zh "D2" "G2"
Add-MirrorHyperlink $zh "A3" "F3"
Add-MirrorHyperlink $zh "D3" "G3"

$zh.Range("H2").Value = "2016-03-21 06:12:49"
$zh.Range("H3").Value = "2016-03-21 06:12:49"

# de-de: row 2 (25c8b8b1...) and row 3 (c582a784...)
Add-MirrorHyperlink $de "A2" "F2"
Add-MirrorHyperlink $de "D2" "G2"
Add-MirrorHyperlink $de "A3" "F3"
Add-MirrorHyperlink $de "D3" "G3"

$de.Range("H2").Value = "2016-03-21 06:12:55"
$de.Range("H3").Value = "2016-03-21 06:12:55"
